# Updates the service-report field values in the filled fault-sheet
# table (date, job number, site name, serial/asset numbers, and the
# fault-description narrative) to the corrected/actual values.

$d = $word.ActiveDocument

function Replace-ExactText($oldText, $newText) {
    # Use a fresh Find each time so each search starts from the top of
    # the document and matches are located unambiguously by their
    # (unique) old value.
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false,
                                $true, 1, $false, "", 0)
    if ($found) {
        # Assign the replacement directly on the matched Range rather
        # than via Find.Execute's Replace argument - this avoids Word's
        # "smart quotes" AutoCorrect mangling a straight apostrophe in
        # the replacement text (e.g. PAK'nSAVE) into a curly one.
        $rng.Text = $newText
    } else {
        Write-Output "WARNING: text not found -> $oldText"
    }
}

Replace-ExactText "23-02-17" "03-10-17"
Replace-ExactText "SV1702200002" "SV1709110022"
Replace-ExactText "New World Rototuna" "Te Awamutu | PAK'nSAVE"
Replace-ExactText "AG22016765E0" "843085"
Replace-ExactText "900015159" "125080"
Replace-ExactText "AK25005042E0" "2189083"
Replace-ExactText "900017086" "110019"
Replace-ExactText "printer destroyed multiple ink ribbons this morning. Replaced with new printer, configured and tested. Working fine." "not opening. Damaged cable. Replaced and tested."
